$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spring")

# ---------------------------------------------------------------------
# The Spring sheet is a repeating 15-row "week" block (header + 7 day
# columns + per-task rows + daily-total row). The last existing block
# lives at rows 154:166. We append a new week block at rows 169:181 by
# copying the previous block and then updating the handful of cells
# that differ (the date row and the task totals that changed).
# ---------------------------------------------------------------------

# --- Body rows (155:166 -> 170:181): plain copy, keeps identical styles ---
$srcBody = $ws.Range("A155:I166")
$dstBody = $ws.Range("A170:I181")
$srcBody.Copy($dstBody)

# --- Header row (154 -> 169): this row is a merged cell (B:H), so copy
# the whole row first (this also registers the B169:H169 merged range),
# then re-apply the clean (unmerged-source) formatting on top so every
# cell in the merge keeps the same single style as the source. ---
$ws.Range("A154:I154").Copy($ws.Range("A169:I169"))
$ws.Range("B154").Copy()
$ws.Range("B169:H169").PasteSpecial(-4122)  # xlPasteFormats

# --- Update the date header row (row 169 label + row 170 dates) ---
$ws.Range("B170").Value = 45389
$ws.Range("C170").Value = 45390
$ws.Range("D170").Value = 45391
$ws.Range("E170").Value = 45026
$ws.Range("F170").Value = 45393
$ws.Range("G170").Value = 45394
$ws.Range("H170").Value = 45395

# --- Update the "Cloud" task row (row 179) hours ---
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = $null
$ws.Range("F179").Value = 3
$ws.Range("G179").Value = 5
$ws.Range("H179").Value = 3
$ws.Range("I179").Value = 16

# --- Update the "Daily Total" row (row 181) ---
$ws.Range("C181").Value = 2
$ws.Range("D181").Value = 0
$ws.Range("F181").Value = 3
$ws.Range("G181").Value = 5
$ws.Range("H181").Value = 3
$ws.Range("I181").Value = 18

# --- Scroll the view down to the newly-added block and select the next
# empty cell below it, mirroring what the user would see after typing
# the new week's data in (matches the saved sheetView/selection). ---
$ws.Activate()
$ws.Range("J181").Select()
$excel.ActiveWindow.ScrollRow = 176

Write-Host "done"
